$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 23.02.2022 21:00"

# Row 10 (EuroOil Opustena): swap current/old price values
$ws.Range("B10").Value = 37.4
$ws.Range("C10").Value = 37.7

# Update the delta price text (keep it stored as text, not a number)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "-0.3"
$ws.Range("D10").Style = "Normal"

# Update the old-date/timestamp text
$ws.Range("E10").Value = "2022-02-23 21:00:29"
